$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 74 values (B..AE) ---
$ws.Cells.Item(74, 2).Value = 7699
$ws.Cells.Item(74, 3).Value = 1036
$ws.Cells.Item(74, 4).Value = 1936
$ws.Cells.Item(74, 5).Value = -901
$ws.Cells.Item(74, 6).Value = -1380
$ws.Cells.Item(74, 7).Value = 109
$ws.Cells.Item(74, 8).Value = -1488
$ws.Cells.Item(74, 9).Value = 466
$ws.Cells.Item(74, 10).Value = 123
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 123
$ws.Cells.Item(74, 13).Value = 3253
$ws.Cells.Item(74, 14).Value = 1924
$ws.Cells.Item(74, 15).Value = -891
$ws.Cells.Item(74, 16).Value = 2220
$ws.Cells.Item(74, 17).Value = -2
$ws.Cells.Item(74, 18).Value = -2
$ws.Cells.Item(74, 19).Value = 4203
$ws.Cells.Item(74, 20).Value = -3683
$ws.Cells.Item(74, 21).Value = 11382
$ws.Cells.Item(74, 22).Value = -7
$ws.Cells.Item(74, 23).Value = 0
$ws.Cells.Item(74, 24).Value = -7
$ws.Cells.Item(74, 25).Value = 1696
$ws.Cells.Item(74, 26).Value = 194
$ws.Cells.Item(74, 27).Value = 1502
$ws.Cells.Item(74, 28).Value = 2963
$ws.Cells.Item(74, 29).Value = 8
$ws.Cells.Item(74, 30).Value = 2954
$ws.Cells.Item(74, 31).Value = 6729

# --- Add new row 75 ---
# Column A holds a text period label ("01-04-2021"); force text format so
# Excel doesn't auto-convert it into a date serial, then restore the
# default "Normal" style so no new cell style gets created.
$ws.Cells.Item(75, 1).NumberFormat = "@"
$ws.Cells.Item(75, 1).Value = "01-04-2021"
$ws.Cells.Item(75, 1).Style = "Normal"

$ws.Cells.Item(75, 2).Value = 3023
$ws.Cells.Item(75, 3).Value = 166
$ws.Cells.Item(75, 4).Value = 2307
$ws.Cells.Item(75, 5).Value = -2142
$ws.Cells.Item(75, 6).Value = 2521
$ws.Cells.Item(75, 7).Value = 1237
$ws.Cells.Item(75, 8).Value = 1284
$ws.Cells.Item(75, 9).Value = 287
$ws.Cells.Item(75, 10).Value = 56
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 56
$ws.Cells.Item(75, 13).Value = -910
$ws.Cells.Item(75, 14).Value = -283
$ws.Cells.Item(75, 15).Value = -480
$ws.Cells.Item(75, 16).Value = -146
$ws.Cells.Item(75, 17).Value = -18
$ws.Cells.Item(75, 18).Value = -18
$ws.Cells.Item(75, 19).Value = 921
$ws.Cells.Item(75, 20).Value = -11
$ws.Cells.Item(75, 21).Value = 3035
$ws.Cells.Item(75, 22).Value = 199
$ws.Cells.Item(75, 23).Value = 0
$ws.Cells.Item(75, 24).Value = 199
$ws.Cells.Item(75, 25).Value = 303
$ws.Cells.Item(75, 26).Value = 50
$ws.Cells.Item(75, 27).Value = 253
$ws.Cells.Item(75, 28).Value = 1203
$ws.Cells.Item(75, 29).Value = 15
$ws.Cells.Item(75, 30).Value = 1188
$ws.Cells.Item(75, 31).Value = 1329
